$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the metric labels in column B (rows 3-21) to reflect the new
# ordering introduced by the "minor bug fixes, swaps added" commit.
$ws.Range("B3").Value  = "pool_balance_sol"
$ws.Range("B4").Value  = "number_of_liquidations"
$ws.Range("B5").Value  = "volume_eth"
$ws.Range("B6").Value  = "treasury_balance"
$ws.Range("B7").Value  = "oi_short"
$ws.Range("B8").Value  = "pool_balance_usdT"
$ws.Range("B9").Value  = "cum_apy_providers"
$ws.Range("B10").Value = "oi_long"
$ws.Range("B11").Value = "min_pnl_traders"
$ws.Range("B12").Value = "fees_collected"
$ws.Range("B13").Value = "pool_balance_btc"
$ws.Range("B14").Value = "pool_balance_eth"
$ws.Range("B15").Value = "cum_pnl_traders"
$ws.Range("B16").Value = "number_of_traders"
$ws.Range("B17").Value = "volume_sol"
$ws.Range("B18").Value = "volume_btc"
$ws.Range("B19").Value = "max_pnl_traders"
$ws.Range("B20").Value = "number_of_liquidity_providers"
$ws.Range("B21").Value = "pool_balance_usdc"
